$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values could otherwise
# be auto-coerced into numbers by Excel (e.g. "1.006" -> 1.006).
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.058.18"
$ws.Range("E2").Value = "  -4.15%  "
$ws.Range("D3").Value = "1.964.33"
$ws.Range("E3").Value = "  -6.19%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "327.45"
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "0.4984"
$ws.Range("E7").Value = "  -6.31%  "
$ws.Range("D8").Value = "0.4205"
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("D9").Value = "52.72"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "0.09209"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "1.097"
$ws.Range("E11").Value = "  -6.98%  "
$ws.Range("D12").Value = "22.92"
$ws.Range("E12").Value = "  -7.34%  "
$ws.Range("D13").Value = "1.993.46"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").Value = "7.849"
$ws.Range("E14").Value = "  -8.72%  "
$ws.Range("D15").Value = "6.438"
$ws.Range("E15").Value = "  -6.48%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "0.00001102"
$ws.Range("E17").Value = "  -5.27%  "
$ws.Range("D18").Value = "91.42"
$ws.Range("E18").Value = "  -10.05%  "
$ws.Range("D19").Value = "0.06715"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "19.16"
$ws.Range("E20").Value = "  -9.42%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "5.969"
$ws.Range("E22").Value = "  -6.07%  "
$ws.Range("D23").Value = "29.090.32"
$ws.Range("E23").Value = "  -4.05%  "
$ws.Range("D24").Value = "12.10"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "2.223.26"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "20.60"
$ws.Range("E27").Value = "  -5.73%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "156.23"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").Value = "6.164"
$ws.Range("E29").Value = "  -10.44%  "
$ws.Range("D30").Value = "2.259"
$ws.Range("E30").Value = "  -9.53%  "
$ws.Range("D31").Value = "126.62"
$ws.Range("E31").Value = "  -5.28%  "
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").Value = "  -8.42%  "
$ws.Range("D33").Value = "0.09833"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").Value = "1.526"
$ws.Range("E34").Value = "  -8.16%  "
$ws.Range("D35").Value = "5.763"
$ws.Range("E35").Value = "  -8.10%  "
$ws.Range("D36").Value = "3.681"
$ws.Range("E36").Value = "  -5.78%  "
$ws.Range("D37").Value = "0.02416"
$ws.Range("E37").Value = "  -8.16%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "0.06327"
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("D40").Value = "8.948"
$ws.Range("E40").Value = "  -11.56%  "
$ws.Range("D41").Value = "0.6440"
$ws.Range("E41").Value = "  -7.77%  "
$ws.Range("D42").Value = "11.41"
$ws.Range("E42").Value = "  -9.74%  "
$ws.Range("D43").Value = "0.1981"
$ws.Range("E43").Value = "  -10.41%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "0.6227"
$ws.Range("E45").Value = "  -8.00%  "
$ws.Range("D46").Value = "13.42"
$ws.Range("E46").Value = "  -6.01%  "
$ws.Range("D47").Value = "2.183"
$ws.Range("E47").Value = "  -8.96%  "
$ws.Range("D48").Value = "1.274"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "3.466"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").Value = "0.00000000330"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").Value = "0.06956"
$ws.Range("E51").Value = "  -4.28%  "
